$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.178.53"
$ws.Range("E2").Value = "  -3.22%  "
Set-TextValue $ws.Range("D3") "3.165.81"
$ws.Range("E3").Value = "  -8.38%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "562.71"
$ws.Range("E5").Value = "  -4.23%  "
Set-TextValue $ws.Range("D6") "170.10"
$ws.Range("E6").Value = "  -3.77%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D8") "0.608"
$ws.Range("E8").Value = "  +0.79%  "
Set-TextValue $ws.Range("D9") "3.163.55"
$ws.Range("E9").Value = "  -8.42%  "
$ws.Range("E10").Value = "  -6.59%  "
$ws.Range("E11").Value = "  -5.14%  "
$ws.Range("E12").Value = "  -4.96%  "
Set-TextValue $ws.Range("D13") "3.712.01"
$ws.Range("E13").Value = "  -8.44%  "
$ws.Range("E14").Value = "  +0.99%  "
Set-TextValue $ws.Range("D15") "27.33"
$ws.Range("E15").Value = "  -7.59%  "
Set-TextValue $ws.Range("D16") "64.129.71"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("E17").Value = "  -5.71%  "
Set-TextValue $ws.Range("D18") "3.161.36"
$ws.Range("E18").Value = "  -8.38%  "
Set-TextValue $ws.Range("D19") "5.73"
$ws.Range("E19").Value = "  -3.93%  "
Set-TextValue $ws.Range("D20") "13.02"
$ws.Range("E20").Value = "  -5.65%  "
Set-TextValue $ws.Range("D21") "353.97"
$ws.Range("E21").Value = "  -5.30%  "
Set-TextValue $ws.Range("D22") "7.22"
$ws.Range("E22").Value = "  -5.35%  "
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E24").Value = "  -5.60%  "
Set-TextValue $ws.Range("D25") "0.504"
$ws.Range("E25").Value = "  -6.52%  "
$ws.Range("E26").Value = "  -5.97%  "
Set-TextValue $ws.Range("D27") "9.58"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("E28").Value = "  -1.67%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -4.00%  "
Set-TextValue $ws.Range("D31") "0.998"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -5.36%  "
Set-TextValue $ws.Range("D33") "22.07"
$ws.Range("E33").Value = "  -6.90%  "
$ws.Range("E35").Value = "  -5.68%  "
$ws.Range("E36").Value = "  -8.02%  "
Set-TextValue $ws.Range("D37") "154.52"
$ws.Range("E37").Value = "  -4.47%  "
Set-TextValue $ws.Range("D38") "0.810"
$ws.Range("E38").Value = "  -8.19%  "
Set-TextValue $ws.Range("D39") "25.82"
$ws.Range("E39").Value = "  -9.28%  "
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("E41").Value = "  -6.08%  "
Set-TextValue $ws.Range("D42") "2.604.11"
$ws.Range("E42").Value = "  -6.15%  "
$ws.Range("E43").Value = "  -7.17%  "
Set-TextValue $ws.Range("D44") "6.01"
$ws.Range("E44").Value = "  -7.00%  "
Set-TextValue $ws.Range("D45") "39.40"
$ws.Range("E45").Value = "  -1.42%  "
Set-TextValue $ws.Range("D46") "0.0657"
$ws.Range("E46").Value = "  -5.10%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "23.92"
$ws.Range("E47").Value = "  -5.23%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D48") "322.41"
$ws.Range("E48").Value = "  -4.79%  "
Set-TextValue $ws.Range("D49") "0.0271"
$ws.Range("E49").Value = "  -7.61%  "
$ws.Range("E50").Value = "  -0.62%  "
